$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to remain plain text (matching the original inlineStr cells) by
# applying a Text number format before the write, then reverting the display
# style afterwards so no stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = '26.303.72'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.690.44'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  +0.12%  '
Set-TextValue $ws.Range("D5") '219.30'
$ws.Range("E5").Value = '  +0.64%  '
Set-TextValue $ws.Range("D6") '0.5267'
$ws.Range("E6").Value = '  +3.72%  '
$ws.Range("E7").Value = '  +0.08%  '
Set-TextValue $ws.Range("D8") '0.2704'
$ws.Range("E8").Value = '  +1.63%  '
Set-TextValue $ws.Range("D9") '0.06436'
$ws.Range("E9").Value = '  +1.38%  '
Set-TextValue $ws.Range("D10") '22.04'
$ws.Range("E10").Value = '  +2.05%  '
Set-TextValue $ws.Range("D11") '0.07462'
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("D12").Value = '1.708.26'
$ws.Range("E12").Value = '  +2.20%  '
Set-TextValue $ws.Range("D13") '4.561'
$ws.Range("E13").Value = '  +0.36%  '
Set-TextValue $ws.Range("D14") '0.5852'
$ws.Range("E14").Value = '  +0.97%  '
Set-TextValue $ws.Range("D15") '0.000008539'
$ws.Range("E15").Value = '  -0.24%  '
Set-TextValue $ws.Range("D16") '64.53'
$ws.Range("E16").Value = '  -0.85%  '
$ws.Range("D17").Value = '26.350.63'
$ws.Range("E17").Value = '  +0.19%  '
Set-TextValue $ws.Range("D18") '4.961'
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("E19").Value = '  +0.17%  '
$ws.Range("E20").Value = '  +0.64%  '
Set-TextValue $ws.Range("D21") '189.69'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  +0.42%  '
Set-TextValue $ws.Range("D23") '1.008'
$ws.Range("E23").Value = '  +0.03%  '
Set-TextValue $ws.Range("D24") '144.78'
$ws.Range("E24").Value = '  +0.67%  '
Set-TextValue $ws.Range("D25") '7.679'
$ws.Range("E25").Value = '  -0.04%  '
Set-TextValue $ws.Range("D26") '0.1233'
$ws.Range("E26").Value = '  +5.37%  '
$ws.Range("E27").Value = '  +1.12%  '
Set-TextValue $ws.Range("D28") '0.06667'
$ws.Range("E28").Value = '  +14.57%  '
Set-TextValue $ws.Range("D29") '1.353'
$ws.Range("E29").Value = '  +5.53%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  +2.10%  '
Set-TextValue $ws.Range("D32") '3.580'
$ws.Range("E32").Value = '  +1.55%  '
Set-TextValue $ws.Range("D33") '1.669'
$ws.Range("E33").Value = '  +1.72%  '
Set-TextValue $ws.Range("D34") '1.029'
$ws.Range("E34").Value = '  +1.99%  '
Set-TextValue $ws.Range("D35") '0.6225'
$ws.Range("E35").Value = '  +4.01%  '
Set-TextValue $ws.Range("D36") '2.394'
$ws.Range("E36").Value = '  +1.27%  '
Set-TextValue $ws.Range("D37") '2.703'
$ws.Range("E37").Value = '  +2.33%  '
Set-TextValue $ws.Range("D38") '6.368'
$ws.Range("E38").Value = '  +5.58%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.109.56'
$ws.Range("E39").Value = '  +2.98%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D40") '0.01625'
$ws.Range("E40").Value = '  +0.78%  '
Set-TextValue $ws.Range("D41") '0.8850'
$ws.Range("E41").Value = '  +2.86%  '
$ws.Range("E42").Value = '  +0.77%  '
Set-TextValue $ws.Range("D43") '100.85'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("D44").Value = '1.838.32'
$ws.Range("E44").Value = '  +0.84%  '
Set-TextValue $ws.Range("D45") '0.00000000113'
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D47") '8.172'
$ws.Range("E47").Value = '  +1.19%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D48") '1.008'
$ws.Range("E48").Value = '  +0.26%  '
Set-TextValue $ws.Range("D49") '0.05268'
$ws.Range("E49").Value = '  +1.55%  '
$ws.Range("E50").Value = '  +0.05%  '
Set-TextValue $ws.Range("D51") '6.059'
$ws.Range("E51").Value = '  +3.37%  '
